$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1723.1428
$ws.Range("I19").Value = 3001
$ws.Range("K19").Value = 3001
$ws.Range("M19").Value = -2826
$ws.Range("H40").Value = 8491.4
$ws.Range("J40").Value = 8190.625
$ws.Range("L40").Value = 8190.625
$ws.Range("N40").Value = -8540.625
$ws.Range("H62").Value = 7333
$ws.Range("I62").Value = 7333
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 7333
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = $null
$ws.Range("N62").Value = -6709
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").Value = $null
$ws.Range("H65").Value = 7333
$ws.Range("I65").Value = 7333
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 36665
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = $null
$ws.Range("N65").Value = -33545
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").Value = $null
$ws.Range("H98").Value = 1126.6875
$ws.Range("I98").Value = 1144.4667
$ws.Range("J98").Value = 860
$ws.Range("K98").Value = 1144.4667
$ws.Range("L98").Value = 860
$ws.Range("M98").Value = 353.5333000000001
$ws.Range("N98").Value = -3856
$ws.Range("H113").Value = 2221.5151
$ws.Range("I113").Value = 2011.9231
$ws.Range("K113").Value = 2011.9231
$ws.Range("M113").Value = 1242.0769
$ws.Range("H122").Value = 1126.6875
$ws.Range("I122").Value = 1144.4667
$ws.Range("J122").Value = 860
$ws.Range("K122").Value = 3433.4001
$ws.Range("L122").Value = 2580
$ws.Range("M122").Value = -983.4000999999998
$ws.Range("N122").Value = -7480
$ws.Range("H137").Value = 1992.75
$ws.Range("I137").Value = 1989.8334
$ws.Range("J137").Value = 2001.5
$ws.Range("K137").Value = 5969.5002
$ws.Range("L137").Value = 6004.5
$ws.Range("M137").Value = -3419.5002
$ws.Range("N137").Value = -11104.5
$ws.Range("H138").Value = 2423.04
$ws.Range("I138").Value = 2229.3333
$ws.Range("K138").Value = 6687.999899999999
$ws.Range("M138").Value = -1547.999899999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6231.346
$ws.Range("I32").Value = 2948.75
$ws.Range("K32").Value = 2948.75
$ws.Range("M32").Value = -2661.75
$ws.Range("H45").Value = 2287.8076
$ws.Range("I45").Value = 1777.2222
$ws.Range("J45").Value = 3436.625
$ws.Range("K45").Value = 1777.2222
$ws.Range("L45").Value = 3436.625
$ws.Range("M45").Value = -1400.2222
$ws.Range("N45").Value = -4190.625
$ws.Range("H122").Value = 2553.7742
$ws.Range("I122").Value = 2191.5
$ws.Range("J122").Value = 4437.6
$ws.Range("K122").Value = 6574.5
$ws.Range("L122").Value = 13312.8
$ws.Range("M122").Value = -4124.5
$ws.Range("N122").Value = -18212.8
$ws.Range("H132").Value = 2003.069
$ws.Range("I132").Value = 1157.3077
$ws.Range("K132").Value = 3471.9231
$ws.Range("M132").Value = -941.9231

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 13105.818
$ws.Range("I20").Value = 3812.75
$ws.Range("J20").Value = 18416.143
$ws.Range("K20").Value = 3812.75
$ws.Range("L20").Value = 18416.143
$ws.Range("M20").Value = -3565.75
$ws.Range("N20").Value = -18910.143
$ws.Range("H86").Value = 2308
$ws.Range("I86").Value = 1997.8
$ws.Range("K86").Value = 1997.8
$ws.Range("M86").Value = -874.8
$ws.Range("H89").Value = 2308
$ws.Range("I89").Value = 1997.8
$ws.Range("K89").Value = 9989
$ws.Range("M89").Value = -4373
$ws.Range("H94").Value = 1075.6
$ws.Range("I94").Value = 913.5
$ws.Range("K94").Value = 913.5
$ws.Range("M94").Value = -462.5
$ws.Range("H134").Value = 5136.091
$ws.Range("I134").Value = 3339
$ws.Range("K134").Value = 10017
$ws.Range("M134").Value = -7482

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3605.9092
$ws.Range("I99").Value = 3477.7
$ws.Range("J99").Value = 4888
$ws.Range("K99").Value = 3477.7
$ws.Range("L99").Value = 4888
$ws.Range("M99").Value = -1979.7
$ws.Range("N99").Value = -7884
$ws.Range("H122").Value = 3661.3635
$ws.Range("J122").Value = 4078.5
$ws.Range("L122").Value = 12235.5
$ws.Range("N122").Value = -17135.5
$ws.Range("H126").Value = 3605.9092
$ws.Range("I126").Value = 3477.7
$ws.Range("J126").Value = 4888
$ws.Range("K126").Value = 10433.1
$ws.Range("L126").Value = 14664
$ws.Range("M126").Value = -7963.099999999999
$ws.Range("N126").Value = -19604

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 999
$ws.Range("J104").Value = 999
$ws.Range("L104").Value = 2997
$ws.Range("N104").Value = -8239
$ws.Range("H106").Value = 18000
$ws.Range("J106").Value = 18000
$ws.Range("L106").Value = 54000
$ws.Range("N106").Value = -55892
$ws.Range("H141").Value = 7577.4443
$ws.Range("I141").Value = 6004.4287
$ws.Range("K141").Value = 18013.2861
$ws.Range("M141").Value = -12833.2861

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2589.0588
$ws.Range("I102").Value = 1078.0769
$ws.Range("K102").Value = 1078.0769
$ws.Range("M102").Value = 543.9231
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = $null
$ws.Range("N109").Value = 0
$ws.Range("H126").Value = 7909.1113
$ws.Range("I126").Value = 7909.1113
$ws.Range("K126").Value = 23727.3339
$ws.Range("M126").Value = -21257.3339
$ws.Range("H132").Value = 5774.852
$ws.Range("I132").Value = 3463.389
$ws.Range("J132").Value = 10397.777
$ws.Range("K132").Value = 10390.167
$ws.Range("L132").Value = 31193.331
$ws.Range("M132").Value = -7860.167000000001
$ws.Range("N132").Value = -36253.331

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7298.5
$ws.Range("I7").Value = 3478
$ws.Range("J7").Value = 13666
$ws.Range("K7").Value = 3478
$ws.Range("L7").Value = 13666
$ws.Range("M7").Value = -3366
$ws.Range("N7").Value = -13890
$ws.Range("H40").Value = 11998.375
$ws.Range("I40").Value = 9996.6
$ws.Range("J40").Value = 15334.667
$ws.Range("K40").Value = 9996.6
$ws.Range("L40").Value = 15334.667
$ws.Range("M40").Value = -9860.6
$ws.Range("N40").Value = -15606.667
$ws.Range("H97").Value = 21666
$ws.Range("J97").Value = 21666
$ws.Range("L97").Value = 21666
$ws.Range("N97").Value = -23648
$ws.Range("H126").Value = 7298.5
$ws.Range("I126").Value = 3478
$ws.Range("J126").Value = 13666
$ws.Range("K126").Value = 10434
$ws.Range("L126").Value = 40998
$ws.Range("M126").Value = -7964
$ws.Range("N126").Value = -45938
$ws.Range("H136").Value = 8004.95
$ws.Range("I136").Value = 5527.2
$ws.Range("J136").Value = 9491.6
$ws.Range("K136").Value = 16581.6
$ws.Range("L136").Value = 28474.8
$ws.Range("M136").Value = -14031.6
$ws.Range("N136").Value = -33574.8

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 3486
$ws.Range("I3").Value = 3486
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 3486
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = $null
$ws.Range("N3").Value = -3372
$ws.Range("H56").Value = 10285
$ws.Range("I56").Value = 10285
$ws.Range("K56").Value = 10285
$ws.Range("M56").Value = -9571
$ws.Range("H96").Value = 10649.625
$ws.Range("J96").Value = 12549.5
$ws.Range("L96").Value = 12549.5
$ws.Range("N96").Value = -15295.5
$ws.Range("H126").Value = 7487.7144
$ws.Range("I126").Value = 7487.7144
$ws.Range("K126").Value = 22463.1432
$ws.Range("M126").Value = -19993.1432
$ws.Range("H132").Value = 2776.0334
$ws.Range("I132").Value = 2612.4119
$ws.Range("K132").Value = 7837.2357
$ws.Range("M132").Value = -5307.2357
$ws.Range("H136").Value = 3743.2307
$ws.Range("I136").Value = 2555.25
$ws.Range("K136").Value = 7665.75
$ws.Range("M136").Value = -5115.75
